# Updated w Mar 31 data
#
# The authoring edit:
#  1. Inserts a new column before the old "AK" (lamesa) column and labels it
#     "imperialbeach", shifting every later city/metric column one to the
#     right.
#  2. Adds a brand-new trailing "other" column.
#  3. Adds a new data row (row 27) for 2020-03-31 (serial 43921).
#  4. Back-fills the "tested" figures for 2020-03-29 / 2020-03-30
#     (rows 25 / 26) which had been left blank.
#  5. Updates the _FilterDatabase defined name so it still spans the full
#     (now one-column-wider) header range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1 & 2: insert the new "imperialbeach" column and add "other" header ---
$ws.Columns("AK:AK").Insert()
$ws.Range("AK1").Value = "imperialbeach"
$ws.Range("BD1").Value = "other"

# --- 4: back-fill "tested" counts that were missing ---
$ws.Range("B25").Value = 687
$ws.Range("B26").Value = 1538

# --- 3: new row 27 (2020-03-31) ---
$ws.Range("A27").Value = 43921
$ws.Range("C27").Value = 734
$ws.Range("D27").Value = 5
$ws.Range("F27").Value = 7
$ws.Range("G27").Value = 1
$ws.Range("H27").Value = 144
$ws.Range("I27").Value = 8
$ws.Range("J27").Value = 174
$ws.Range("K27").Value = 17
$ws.Range("L27").Value = 130
$ws.Range("M27").Value = 18
$ws.Range("N27").Value = 116
$ws.Range("O27").Value = 27
$ws.Range("P27").Value = 74
$ws.Range("Q27").Value = 20
$ws.Range("R27").Value = 49
$ws.Range("S27").Value = 22
$ws.Range("T27").Value = 34
$ws.Range("U27").Value = 23
$ws.Range("V27").Value = 1
$ws.Range("W27").Value = 0
$ws.Range("X27").Value = 317
$ws.Range("Y27").Value = 414
$ws.Range("Z27").Value = 3
$ws.Range("AA27").Value = 136
$ws.Range("AB27").Value = 56
$ws.Range("AC27").Value = 9
$ws.Range("AD27").Value = 27
$ws.Range("AE27").Value = 38
$ws.Range("AF27").Value = 1
$ws.Range("AG27").Value = 5
$ws.Range("AH27").Value = 36
$ws.Range("AI27").Value = 20
$ws.Range("AJ27").Value = 11
$ws.Range("AL27").Value = 8
$ws.Range("AM27").Value = 5
$ws.Range("AN27").Value = 10
$ws.Range("AO27").Value = 17
$ws.Range("AP27").Value = 7
$ws.Range("AQ27").Value = 422
$ws.Range("AR27").Value = 5
$ws.Range("AS27").Value = 5
$ws.Range("AT27").Value = 3
$ws.Range("AU27").Value = 8
$ws.Range("AV27").Value = 3
$ws.Range("AW27").Value = 4
$ws.Range("AY27").Value = 5
$ws.Range("AZ27").Value = 4
$ws.Range("BA27").Value = 11
$ws.Range("BC27").Value = 17
$ws.Range("BD27").Value = 26

# --- 5: fix up the hidden _FilterDatabase name (AZ -> BA since a column
#         was inserted before it) ---
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Sheet1!_FilterDatabase") {
        $n.RefersTo = "=Sheet1!`$A`$1:`$BA`$17"
    }
}

# --- keep the active selection in sync with the author's final click ---
$ws.Range("Y27").Select()
